$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "thanhtoan"
$ws.Range("D1").Value = "loaixe"
$ws.Range("E1").Value = "ngaydat"
$ws.Range("F1").Value = "ngaylay"
$ws.Range("G1").Value = "ngaytra"

$ws.Range("D3").Select()
